$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 3.1851815
$ws.Range("N2").Value = 6.370363
$ws.Range("O2").Value = 0.4406530230187619
$ws.Range("P2").Value = 0.3851702893788179
$ws.Range("Q2").Value = 0.38642144180775
$ws.Range("R2").Value = 1.545685767231
$ws.Range("S2").Value = 0.4406530230187619
$ws.Range("T2").Value = 0.3851702893788179

# Row 3
$ws.Range("O3").Value = 0.2827048402157753
$ws.Range("P3").Value = 0.3706641033643825
$ws.Range("S3").Value = 0.2827048402157753
$ws.Range("T3").Value = 0.3706641033643825

# Row 4
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.003928
$ws.Range("N4").Value = 0.011784
$ws.Range("O4").Value = 0.0005434180358066555
$ws.Range("P4").Value = 0.0007124942001013113
$ws.Range("Q4").Value = 0.000476539068
$ws.Range("R4").Value = 0.002859234408
$ws.Range("S4").Value = 0.0005434180358066555
$ws.Range("T4").Value = 0.0007124942001013113

# Row 5
$ws.Range("M5").Value = 1.9606995
$ws.Range("N5").Value = 3.921399
$ws.Range("O5").Value = 0.2712524111754306
$ws.Range("P5").Value = 0.2370989514411984
$ws.Range("Q5").Value = 0.23786912229075
$ws.Range("R5").Value = 0.951476489163
$ws.Range("S5").Value = 0.2712524111754306
$ws.Range("T5").Value = 0.2370989514411984

# Row 6
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.02851766666666667
$ws.Range("N6").Value = 0.085553
$ws.Range("O6").Value = 0.003945268433245655
$ws.Range("P6").Value = 0.005172778029639129
$ws.Range("Q6").Value = 0.0034597205435
$ws.Range("R6").Value = 0.020758323261
$ws.Range("S6").Value = 0.003945268433245655
$ws.Range("T6").Value = 0.005172778029639129

# Row 7
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.006513000000000001
$ws.Range("N7").Value = 0.019539
$ws.Range("O7").Value = 0.0009010391209798237
$ws.Range("P7").Value = 0.001181383585860448
$ws.Range("Q7").Value = 0.0007901473905
$ws.Range("R7").Value = 0.004740884343
$ws.Range("S7").Value = 0.0009010391209798237
$ws.Range("T7").Value = 0.001181383585860448
